$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Replace the 4M12/1M divider (R1-R4 @ 1M, R10 @ 4M12) with a
# 3M9/953k divider (R1-R3 @ 1M, R4 @ 953k, R10 @ 3M9), and fix the
# IC1 part number (ATTINY84A-SSU -> ATTINY44A-SSU), keeping the
# workbook within JLCPCB's basic parts library.
# ------------------------------------------------------------------

# 1) Insert a new row at 12 so R4 can become its own BOM line right
#    after the remaining R1/R2/R3 @ 1M line.
$ws.Rows.Item(12).Insert()

# R1, R2, R3, R4 @ 1M  ->  R1, R2, R3 @ 1M
$ws.Range("A11").Value = "R1, R2, R3"

# New line: R4 @ 953k
$ws.Range("A12").Value = "R4"
$ws.Range("B12").Value = "953k"
$ws.Range("C12").Value = "R0603"
$ws.Range("D12").Value = "C246846"
$ws.Range("E12").Value = "e"

# 2) The old R10 @ 4M12 line (originally row 21) is now at row 22
#    after the insert above. Remove it - R10 will be re-inserted
#    below with its new value.
$ws.Rows.Item(22).Delete()

# 3) Insert a new row at 20 (just before IC3) for R10 @ 3M9.
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "R10"
$ws.Range("B20").Value = "3M9"
$ws.Range("C20").Value = "R0603"
$ws.Range("D20").Value = "C23019"
$ws.Range("E20").Value = "e"

# 4) Fix the IC1 comment (now at row 24).
$ws.Range("B24").Value = "ATTINY44A-SSU"

# 5) Resize the "Daten" table so it covers the extra row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E33"))
